$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Updated path to NVIS Extant file (new date 20240801, was 20240709)
$ws.Range("B4").Value = "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\processing\NEAP_intermediate\NVIS_IUCNGET_DK_20240801.tif"

# Match the author's final selection in the saved workbook view
$ws.Range("B5").Select() | Out-Null
